$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.489400148391724
$ws.Range("B1").Value = 3.613487958908081
$ws.Range("C1").Value = 2.491650104522705
$ws.Range("D1").Value = 1.298937797546387
$ws.Range("E1").Value = 0.7618328332901001
